$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing the existing rows 4-26 down to 5-27.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly data point.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44630
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 100112045
$ws.Range("G4").Value = "Zapallo"
$ws.Range("H4").Value = "Camote"
$ws.Range("I4").Value = "1a (cosecha)"
$ws.Range("J4").Value = 1200
$ws.Range("K4").Value = 450
$ws.Range("L4").Value = 480
$ws.Range("M4").Value = 465
$ws.Range("N4").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 465
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Hortaliza"
